$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 5
$ws.Range("I3").Value = 5
$ws.Range("I4").Value = 5
$ws.Range("I5").Value = 5

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 6
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 5
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = -1
$ws.Range("G6").Value = 5
$ws.Range("H6").Value = 21
$ws.Range("I6").Value = 5
$ws.Range("J6").Value = "train_dim2_1"

$ws.Range("I7").Select()
